$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are plain numeric-looking strings (prices using "." as
# thousands separators). Force them to stay text, matching the workbook's
# original inlineStr storage, instead of letting Excel auto-convert to Number.

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "28.450.63"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  -2.81%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.771.41"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  -3.30%  "

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.002"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  -0.23%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "229.05"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -2.28%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.5858"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -2.41%  "

$ws.Range("E7").Value = "  -0.24%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.2735"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -0.83%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "23.13"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -0.64%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.06673"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -4.41%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.07533"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -1.10%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "1.768.64"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -3.59%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "4.745"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -0.17%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "0.6054"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -3.45%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "2.012.92"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -3.17%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "74.48"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -5.07%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "0.000008583"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -11.20%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "28.432.89"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -1.09%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "5.351"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -5.19%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "1.002"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -0.22%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "206.18"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -6.53%  "

$ws.Range("E22").Value = "  -1.79%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "6.718"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -2.20%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "1.004"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -0.32%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "151.21"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -3.30%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "8.097"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +1.72%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "0.1245"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -3.15%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "16.22"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -1.91%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "1.406"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -3.27%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "0.06144"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -3.92%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "1.409"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -2.14%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "3.749"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -2.34%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "3.746"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -0.21%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "1.669"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -2.98%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "1.040"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -4.55%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.6342"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -1.73%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "2.500"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -1.61%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "2.671"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -2.84%  "

$ws.Range("E39").Value = "  -5.29%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "6.306"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -4.20%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "1.132.16"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -2.47%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.8707"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -2.42%  "

$ws.Range("E43").Value = "  +0.18%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "99.82"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -0.59%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "1.926.46"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -2.99%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "59.45"
$cell.Style = "Normal"

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.00000000109"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -3.37%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "1.570"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -1.05%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "8.328"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -1.56%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.05408"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -2.51%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.4462"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -2.08%  "
